$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 4) mirroring rows 2-3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1.1388310185185186
$ws.Range("B4").NumberFormat = $ws.Range("B3").NumberFormat
$ws.Range("C4").Value = "H2O: Just Add Water Season 2 (Audiovisual, English, Familiar):20; Berserk (Text with visuals, Japanese, Familiar):23; Ijiranaide, Nagatoro-san(Text with visuals, Japanese, Re-watch):28;"
$ws.Range("D4").Value = "Watched children's shows I'm familiar with and read simple manga."

# Match the active selection Excel leaves after data entry
$ws.Range("C4").Select()
